$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "IVOS" team is being dropped from the upper-triangle comparison: remove
# both its column (Q) and its row (17) from the matrix.
$ws.Columns("Q").Delete()
$ws.Rows("17").Delete()

# Clean up / normalize a couple of team name labels (both the column header in
# row 1 and the row label in column A for each team).
$ws.Range("A4").Value = "VIREO"
$ws.Range("D1").Value = "VIREO"

$ws.Range("A5").Value = "SOMHunter"
$ws.Range("E1").Value = "SOMHunter"

# Shrink the conditional formatting (color scale) range so it keeps covering
# exactly the data area after the row/column removal.
$rng = $ws.Range("B2:R18")
$cfs = $rng.FormatConditions
$cf1 = $cfs.Item(1)
$cf1.ModifyAppliesToRange($ws.Range("B2:Q17"))
$cf1.Priority = 3

# Leave the selection where the author ended up after the cleanup pass.
$ws.Range("I9").Select() | Out-Null
